{"js": "// Update the date line and the twenty-five two-digit multiplication\n// problems scattered through the table, matching the source docx's\n// before/after text pairs 1:1 (every old value is unique in the doc).\nconst replacements = [\n  [\"2025-08-06 Wednesday\", \"2025-08-07 Thursday\"],\n  [\"58\u00d739=\", \"46\u00d715=\"],\n  [\"50\u00d757=\", \"62\u00d724=\"],\n  [\"16\u00d749=\", \"83\u00d767=\"],\n  [\"14\u00d753=\", \"40\u00d793=\"],\n  [\"12\u00d760=\", \"95\u00d748=\"],\n  [\"44\u00d733=\", \"19\u00d725=\"],\n  [\"39\u00d781=\", \"61\u00d784=\"],\n  [\"84\u00d799=\", \"19\u00d792=\"],\n  [\"91\u00d724=\", \"21\u00d789=\"],\n  [\"31\u00d784=\", \"74\u00d733=\"],\n  [\"37\u00d779=\", \"17\u00d784=\"],\n  [\"13\u00d719=\", \"72\u00d762=\"],\n  [\"57\u00d780=\", \"68\u00d720=\"],\n  [\"92\u00d729=\", \"31\u00d752=\"],\n  [\"61\u00d763=\", \"33\u00d778=\"],\n  [\"68\u00d735=\", \"50\u00d729=\"],\n  [\"70\u00d732=\", \"29\u00d742=\"],\n  [\"76\u00d797=\", \"32\u00d759=\"],\n  [\"88\u00d741=\", \"40\u00d721=\"],\n  [\"65\u00d767=\", \"53\u00d714=\"],\n  [\"91\u00d725=\", \"92\u00d790=\"],\n  [\"23\u00d742=\", \"80\u00d775=\"],\n  [\"76\u00d788=\", \"34\u00d728=\"],\n  [\"56\u00d720=\", \"69\u00d714=\"],\n  [\"41\u00d796=\", \"37\u00d777=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the twenty-five two-digit multiplication\n# problems scattered through the table. Each \"old\" value is unique in\n# the document, so a simple Find/Replace per pair is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-08-06 Wednesday\", \"2025-08-07 Thursday\"),\n    @(\"58\u00d739=\", \"46\u00d715=\"),\n    @(\"50\u00d757=\", \"62\u00d724=\"),\n    @(\"16\u00d749=\", \"83\u00d767=\"),\n    @(\"14\u00d753=\", \"40\u00d793=\"),\n    @(\"12\u00d760=\", \"95\u00d748=\"),\n    @(\"44\u00d733=\", \"19\u00d725=\"),\n    @(\"39\u00d781=\", \"61\u00d784=\"),\n    @(\"84\u00d799=\", \"19\u00d792=\"),\n    @(\"91\u00d724=\", \"21\u00d789=\"),\n    @(\"31\u00d784=\", \"74\u00d733=\"),\n    @(\"37\u00d779=\", \"17\u00d784=\"),\n    @(\"13\u00d719=\", \"72\u00d762=\"),\n    @(\"57\u00d780=\", \"68\u00d720=\"),\n    @(\"92\u00d729=\", \"31\u00d752=\"),\n    @(\"61\u00d763=\", \"33\u00d778=\"),\n    @(\"68\u00d735=\", \"50\u00d729=\"),\n    @(\"70\u00d732=\", \"29\u00d742=\"),\n    @(\"76\u00d797=\", \"32\u00d759=\"),\n    @(\"88\u00d741=\", \"40\u00d721=\"),\n    @(\"65\u00d767=\", \"53\u00d714=\"),\n    @(\"91\u00d725=\", \"92\u00d790=\"),\n    @(\"23\u00d742=\", \"80\u00d775=\"),\n    @(\"76\u00d788=\", \"34\u00d728=\"),\n    @(\"56\u00d720=\", \"69\u00d714=\"),\n    @(\"41\u00d796=\", \"37\u00d777=\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $find.Execute(\n        $old,\n        $true,\n        $false,\n        $false,\n        $false,\n        $false,\n        $true,\n        1,\n        $false,\n        $new,\n        2\n    ) | Out-Null\n}\n\n$d.Save()\n"}
